$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Switch the workbook's default/body font to Arial
$wb.Styles.Item(1).Font.Name = "Arial"

# Fill in the new "tasks / artifacts" columns (F-I) for rows 1-5
$ws.Range("F4").Value = "User journey map"
$ws.Range("F2").Value = "Affinity map"
$ws.Range("F3").Value = "Persona , Affinity map"
$ws.Range("G2").Value = "Site map"
$ws.Range("F1").Value = "User journey map, Site map"
$ws.Range("F5").Value = "Site map,user flow"
$ws.Range("H5").Value = "Wireframe"

# Correct the wording of an existing note
$ws.Range("E3").Value = "user interview,Competitor analysis"

$ws.Range("I1").Value = "User flow,wireframe"
$ws.Range("H3").Value = "Site map"

# Remove the last team member row (Aya Mohamed), which is no longer present
$ws.Range("A6").ClearContents()

# Match the author's final selection/cursor position
$ws.Range("J6").Select() | Out-Null
